# Fix a typo in the "responsibility_options" column (E): the field-option
# token was written as "IM_operators" (underscore) in several rows while the
# canonical option elsewhere in the sheet uses a hyphen ("IM-operators").
# Normalize every occurrence in column E to use the hyphenated form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -ne $null -and $val.ToString().Contains("IM_operators")) {
        $cell.Value2 = $val.ToString().Replace("IM_operators", "IM-operators")
    }
}
